# ==========================================================================
# Add Q3-2022 holdings data:
#   1. Insert a new "2022-Q3" quarter row at the top of the "总计" summary
#      sheet, shifting the six existing quarters down by one row.
#   2. Insert a brand-new "2022-Q3" worksheet (positioned right after "总计"
#      and before "2022-Q2") containing the full fund-holdings breakdown.
# ==========================================================================

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Step 1: update the "总计" summary sheet ---------------------------------

# Capture the six existing quarter rows (2022-Q2 .. 2021-Q1) before shifting.
$existing = @()
for ($r = 2; $r -le 7; $r++) {
    $existing += ,@($ws1.Cells.Item($r,2).Value(), $ws1.Cells.Item($r,3).Value(), $ws1.Cells.Item($r,4).Value())
}

# Extend the bordered index-column style down into the new last row (row 8)
# by copying the format of the previous last row (row 7).
$ws1.Cells.Item(7,1).Copy($ws1.Cells.Item(8,1))

# Shift the six existing quarters down by one row (old row r -> new row r+1).
for ($i = 5; $i -ge 0; $i--) {
    $newRow = $i + 3
    $ws1.Cells.Item($newRow,1).Value = $i + 1
    $ws1.Cells.Item($newRow,2).Value = $existing[$i][0]
    $ws1.Cells.Item($newRow,3).Value = $existing[$i][1]
    $ws1.Cells.Item($newRow,4).Value = $existing[$i][2]
}

# Write the new 2022-Q3 row at row 2.
$ws1.Cells.Item(2,1).Value = 0
$ws1.Cells.Item(2,2).Value = "2022-Q3"
$ws1.Cells.Item(2,3).Value = 29
$ws1.Cells.Item(2,4).Value = 8.36

# --- Step 2: insert the new "2022-Q3" worksheet ------------------------------

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "2022-Q3"

$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($c = 2; $c -le 8; $c++) {
    # Copy the header style from the "总计" sheet's own header cell, then set the text.
    $ws1.Range("B1").Copy($ws2.Cells.Item(1,$c))
    $ws2.Cells.Item(1,$c).Value = $headers[$c - 2]
}

$q3Data = @(
    @("012930","中庚价值先锋股票","47.83","94.71","8.52","4.0751",1),
    @("920003","中金新锐股票A","17.72","89.26","4.35","0.7708",4),
    @("501078","广发科创主题灵活配置混合（LOF）","12.77","93.49","3.11","0.3971",8),
    @("001487","宝盈优势产业灵活配置混合A","10.11","91.85","3.72","0.3761",1),
    @("001128","宝盈新兴产业灵活配置混合A","9.00","91.66","4.04","0.3636",5),
    @("013895","宝盈成长精选混合A","8.51","90.34","3.93","0.3344",2),
    @("001877","宝盈国家安全沪港深股票A","8.74","91.81","3.71","0.3243",2),
    @("162720","广发创业板两年定期开放混合","6.33","94.20","4.07","0.2576",7),
    @("000586","景顺长城中小创精选股票","2.21","93.50","9.03","0.1996",1),
    @("001075","宝盈转型动力灵活配置混合A","4.35","91.90","3.96","0.1723",4),
    @("001103","前海开源工业革命4.0灵活配置混合","3.27","61.94","4.46","0.1458",4),
    @("920923","中金新锐股票C","3.32","89.26","4.35","0.1444",4),
    @("012771","宝盈优势产业灵活配置混合C","3.62","91.85","3.72","0.1347",1),
    @("013896","宝盈成长精选混合C","2.93","90.34","3.93","0.1151",2),
    @("014339","长江智能制造混合A","2.42","75.36","4.03","0.0975",3),
    @("160642","鹏华增瑞灵活配置混合（LOF）","2.05","90.96","4.17","0.0855",7),
    @("009353","浙商科技创新一个月滚动持有混合A","2.06","91.05","4.10","0.0845",10),
    @("009354","浙商科技创新一个月滚动持有混合C","1.45","91.05","4.10","0.0594",10),
    @("012815","宝盈新兴产业灵活配置混合C","1.31","91.66","4.04","0.0529",5),
    @("260115","景顺长城中小盘混合","0.92","92.87","5.59","0.0514",3),
    @("010706","景顺长城景骊成长混合","0.61","92.93","6.42","0.0392",6),
    @("519644","银河智联主题灵活配置混合","0.90","93.05","3.98","0.0358",10),
    @("011214","招商惠润一年定期开放混合（MOM）A","0.48","57.54","2.50","0.0120",5),
    @("015389","宝盈转型动力灵活配置混合C","0.24","91.90","3.96","0.0095",4),
    @("009128","明亚价值长青混合A","0.38","57.73","2.51","0.0095",10),
    @("013613","宝盈国家安全沪港深股票C","0.23","91.81","3.71","0.0085",2),
    @("014340","长江智能制造混合C","0.11","75.36","4.03","0.0044",3),
    @("011215","招商惠润一年定期开放混合（MOM）C","0.06","57.54","2.50","0.0015",5),
    @("009129","明亚价值长青混合C","0.00","57.73","2.51","0",10)
)

# Force columns B (fund code, to preserve leading zeros) and D:G (percentages /
# figures stored as text in this workbook) to Text format before writing, so
# Excel does not silently convert the numeric-looking strings to numbers.
$ws2.Range("B2:B30").NumberFormat = "@"
$ws2.Range("D2:G30").NumberFormat = "@"

for ($i = 0; $i -lt $q3Data.Count; $i++) {
    $r = $i + 2
    $row = $q3Data[$i]
    $ws2.Cells.Item($r,1).Value = $i
    $ws2.Cells.Item($r,2).Value = $row[0]
    $ws2.Cells.Item($r,3).Value = $row[1]
    $ws2.Cells.Item($r,4).Value = $row[2]
    $ws2.Cells.Item($r,5).Value = $row[3]
    $ws2.Cells.Item($r,6).Value = $row[4]
    $ws2.Cells.Item($r,7).Value = $row[5]
    $ws2.Cells.Item($r,8).Value = $row[6]
}

# Row 30's "持有市值" (G30) is stored as a genuine numeric 0 in the source data
# (every other row in that column is textual), so fix it up to match.
$ws2.Range("G30").NumberFormat = "General"
$ws2.Cells.Item(30,7).Value = 0

